# Expand the abbreviation "FA" in "FA_on_Flight(" to "Flight_Attendants"
# so the run reads "Flight_Attendants_on_Flight(".
#
# The target OOXML keeps the single original run ("FA_on_Flight(") but
# splits it into four sibling runs - "F", "light_Attendant", "s",
# "_on_Flight( - all sharing the exact same run properties
# (rFonts ascii/hAnsi/cs = David). A plain Range.Text (or Find/Replace)
# assignment always re-merges same-formatted neighbouring text back into
# a single run, so instead we clone an already-isolated, identically
# formatted run ("F") via FormattedText to mint new, independent run
# boundaries, then rename the clones in place (also through
# FormattedText, which does not trigger the merge-on-write behaviour).

$d = $word.ActiveDocument

# Locate the run we need to edit.
$match = $d.Content
$found = $match.Find.Execute("FA_on_Flight(", $true, $false, $false, $false, `
                              $false, $true, 1, $false, "", 0)
$start = $match.Start

# Template run: the leading "F" - already its own run with the rFonts
# (David/David/David) formatting shared by the whole match.
$template = $d.Range($start, $start + 1)

# --- Split off "A" (2nd char) into its own run, then rename it to
#     "light_Attendant" -----------------------------------------------
$aRun = $d.Range($start + 1, $start + 2)
$aRun.FormattedText = $template.FormattedText

$newRun = $d.Range($start + 1, $start + 2)
$ft = $newRun.FormattedText
$ft.Text = "light_Attendant"
$renamedRun = $d.Range($start + 1, $start + 1 + $ft.Text.Length)
$renamedRun.FormattedText = $ft

# --- Insert a new run for the trailing "s" right after
#     "light_Attendant" (and before "_on_Flight(") -----------------------
$insertAt = $start + 1 + $ft.Text.Length
$sSpot = $d.Range($insertAt, $insertAt)
$sSpot.FormattedText = $template.FormattedText

$newSRun = $d.Range($insertAt, $insertAt + 1)
$sFt = $newSRun.FormattedText
$sFt.Text = "s"
$renamedSRun = $d.Range($insertAt, $insertAt + $sFt.Text.Length)
$renamedSRun.FormattedText = $sFt
